# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp
# - Update Estados Unidos (row 4) totals
# - España overtakes Chile in the ranking (rows 12/13 swap countries + data)
# - Refresh a handful of other country rows (Libia, Bahamas, Antigua y Barbuda)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Agosto de 2020 a las 03:20"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 5746272
$ws.Range("C4").Value = 45341
$ws.Range("D4").Value = 3095365
$ws.Range("E4").Value = 2473483
$ws.Range("G4").Value = 1090
$ws.Range("H4").Value = 177424

# Row 12 now holds España (overtook Chile)
$ws.Range("A12").Value = "España"
$ws.Range("B12").Value = 404229
$ws.Range("C12").Value = 3349
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("G12").Value = 16
$ws.Range("H12").Value = 28813

# Row 13 now holds Chile
$ws.Range("A13").Value = "Chile"
$ws.Range("B13").Value = 391849
$ws.Range("C13").Value = 1812
$ws.Range("D13").Value = 366063
$ws.Range("E13").Value = 15115
$ws.Range("G13").Value = 93
$ws.Range("H13").Value = 10671

# Libia (row 91)
$ws.Range("B91").Value = 9707
$ws.Range("C91").Value = 244
$ws.Range("D91").Value = 1047
$ws.Range("E91").Value = 8487
$ws.Range("G91").Value = 4
$ws.Range("H91").Value = 173

# Bahamas (row 143)
$ws.Range("B143").Value = 1610
$ws.Range("C143").Value = 79
$ws.Range("D143").Value = 211
$ws.Range("E143").Value = 1376
$ws.Range("G143").Value = 1
$ws.Range("H143").Value = 23

# Antigua y Barbuda (row 195)
$ws.Range("D195").Value = 89
$ws.Range("E195").Value = 2
